$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.884.22'
$ws.Range("E2").Value = '  +2.65%  '

$ws.Range("D3").Value = '3.725.78'
$ws.Range("E3").Value = '  +6.21%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '420.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").Value = '3.717.31'
$ws.Range("E7").Value = '  +6.20%  '

$ws.Range("E8").Value = '  -0.33%  '

$ws.Range("E9").Value = '  +0.03%  '

$ws.Range("E10").Value = '  -0.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.184'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +13.95%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000407'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +55.95%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '43.08'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.24%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.17%  '

$ws.Range("D15").Value = '4.292.96'
$ws.Range("E15").Value = '  +5.77%  '

$ws.Range("E16").Value = '  -0.84%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.83'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.94%  '

$ws.Range("D18").Value = '3.721.58'
$ws.Range("E18").Value = '  +5.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.31'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.88%  '

$ws.Range("E20").Value = '  +4.34%  '

$ws.Range("D21").Value = '66.931.86'
$ws.Range("E21").Value = '  +2.89%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '449.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.72%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '16.69'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +25.52%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '90.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.07%  '

$ws.Range("E25").Value = '  -1.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '38.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +12.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.30%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.35'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.04'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.76%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.80'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.31%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.125'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +9.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.76'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.43%  '

$ws.Range("E33").Value = '  -3.18%  '

$ws.Range("E34").Value = '  +1.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '42.05'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.23%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.33'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.41%  '

$ws.Range("E37").Value = '  +0.01%  '

$ws.Range("E38").Value = '  -0.47%  '

$ws.Range("D39").Value = '0.0₃0753'
$ws.Range("E39").Value = '  +5.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.11'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +34.03%  '

$ws.Range("E41").Value = '  +1.71%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '28.79'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +33.88%  '

$ws.Range("E43").Value = '  -0.11%  '

$ws.Range("E44").Value = '  +4.83%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.29'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +34.90%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '148.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.09%  '

$ws.Range("E47").Value = '  +6.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.97%  '

$ws.Range("E49").Value = '  -3.58%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.67%  '

$ws.Range("E51").Value = '  -1.10%  '
